# Updates current market-price-derived figures on each job sheet's Leve table.
# Values below are taken from the authoritative post-edit snapshot; only the
# cells that actually changed are touched (comments show the Leve name for context).
$wb = $excel.ActiveWorkbook

# ==== ALC ====
$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 1898.9143
$ws.Range("J17").Value = 1898.9143
$ws.Range("L17").Value = 5696.742899999999
$ws.Range("N17").Value = -6032.742899999999
# Row 32: Automata for the People
$ws.Range("H32").Value = 1216.4286
$ws.Range("J32").Value = 1023
$ws.Range("L32").Value = 1023
$ws.Range("N32").Value = -1675
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 15040.6
$ws.Range("I62").Value = 17139.385
$ws.Range("K62").Value = 17139.385
$ws.Range("M62").Value = -16515.385
# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 15040.6
$ws.Range("I65").Value = 17139.385
$ws.Range("K65").Value = 85696.92499999999
$ws.Range("M65").Value = -82576.92499999999
# Row 69: Steeling the Knife, Steeling the Mind
$ws.Range("H69").Value = 8679.15
$ws.Range("I69").Value = 5404.3335
$ws.Range("K69").Value = 16213.0005
$ws.Range("M69").Value = -15339.0005
# Row 70: Consecrating Congregation
$ws.Range("H70").Value = 13503.833
$ws.Range("I70").Value = 1729.8
$ws.Range("J70").Value = 21913.857
$ws.Range("K70").Value = 5189.4
$ws.Range("L70").Value = 65741.571
$ws.Range("M70").Value = -4919.4
$ws.Range("N70").Value = -66281.571
# Row 72: Surgical Substitution (L)
$ws.Range("H72").Value = 8679.15
$ws.Range("I72").Value = 5404.3335
$ws.Range("K72").Value = 48639.0015
$ws.Range("M72").Value = -44271.0015
# Row 73: Curbing the Contagion (L)
$ws.Range("H73").Value = 13503.833
$ws.Range("I73").Value = 1729.8
$ws.Range("J73").Value = 21913.857
$ws.Range("K73").Value = 5189.4
$ws.Range("L73").Value = 65741.571
$ws.Range("M73").Value = -4253.4
$ws.Range("N73").Value = -67613.571
# Row 93: Spellbound
$ws.Range("H93").Value = 35546
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 2137.8235
$ws.Range("I137").Value = 1229.625
$ws.Range("J137").Value = 2417.2693
$ws.Range("K137").Value = 3688.875
$ws.Range("L137").Value = 7251.8079
$ws.Range("M137").Value = -1138.875
$ws.Range("N137").Value = -12351.8079
# Row 138: All-night Crafting
$ws.Range("H138").Value = 2147.2
$ws.Range("I138").Value = 1541.0333
$ws.Range("K138").Value = 4623.0999
$ws.Range("M138").Value = 516.9000999999998
# Row 141: Remedy for Reason
$ws.Range("H141").Value = 2461.8
$ws.Range("I141").Value = 1999.6666
$ws.Range("J141").Value = 3155
$ws.Range("K141").Value = 5998.9998
$ws.Range("L141").Value = 9465
$ws.Range("M141").Value = -818.9997999999996
$ws.Range("N141").Value = -19825

# ==== ARM ====
$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth
$ws.Range("H5").Value = 160.9375
$ws.Range("I5").Value = 100.36364
$ws.Range("K5").Value = 100.36364
$ws.Range("M5").Value = 11.63636
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 2168.6445
$ws.Range("I74").Value = 1417.3889
$ws.Range("K74").Value = 1417.3889
$ws.Range("M74").Value = -543.3888999999999
# Row 75: Someone Put Dung in My Helmet
$ws.Range("H75").Value = 65813.14
$ws.Range("J75").Value = 80138.39999999999
$ws.Range("L75").Value = 80138.39999999999
$ws.Range("N75").Value = -81886.39999999999
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 2168.6445
$ws.Range("I77").Value = 1417.3889
$ws.Range("K77").Value = 7086.9445
$ws.Range("M77").Value = -2718.9445
# Row 78: Rage against the Scream (L)
$ws.Range("H78").Value = 65813.14
$ws.Range("J78").Value = 80138.39999999999
$ws.Range("L78").Value = 240415.2
$ws.Range("N78").Value = -249151.2
# Row 101: Art Imitates Life
$ws.Range("H101").Value = 94365.60000000001
$ws.Range("J101").Value = 94365.60000000001
$ws.Range("L101").Value = 94365.60000000001
$ws.Range("N101").Value = -100855.6
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3199.81
$ws.Range("I132").Value = 1850.9828
$ws.Range("K132").Value = 5552.9484
$ws.Range("M132").Value = -3022.9484

# ==== BSM ====
$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences
$ws.Range("H4").Value = 160.9375
$ws.Range("I4").Value = 100.36364
$ws.Range("K4").Value = 100.36364
$ws.Range("M4").Value = 14.63636
# Row 96: Hammer Time
$ws.Range("H96").Value = 26792.75
$ws.Range("I96").Value = 10723.667
$ws.Range("K96").Value = 10723.667
$ws.Range("M96").Value = -7977.666999999999
# Row 138: Bladewinner
$ws.Range("H138").Value = 99965.664
$ws.Range("J138").Value = 99965.664
$ws.Range("L138").Value = 99965.664
$ws.Range("N138").Value = -110245.664

# ==== CRP ====
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 4631.896
$ws.Range("J31").Value = 5579.0347
$ws.Range("L31").Value = 5579.0347
$ws.Range("N31").Value = -6169.0347
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 4631.896
$ws.Range("J34").Value = 5579.0347
$ws.Range("L34").Value = 5579.0347
$ws.Range("N34").Value = -5983.0347
# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 5590.5
$ws.Range("I62").Value = 3839.8
$ws.Range("K62").Value = 3839.8
$ws.Range("M62").Value = -3215.8
# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 5590.5
$ws.Range("I65").Value = 3839.8
$ws.Range("K65").Value = 19199
$ws.Range("M65").Value = -16079
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1934.5714
$ws.Range("I134").Value = 1915.3914
$ws.Range("K134").Value = 5746.174199999999
$ws.Range("M134").Value = -3211.174199999999

# ==== CUL ====
$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up
$ws.Range("H12").Value = 28.5
$ws.Range("J12").Value = 39
$ws.Range("L12").Value = 117
$ws.Range("N12").Value = -463
# Row 23: Sweet Smell of Success
$ws.Range("H23").Value = 624.6667
$ws.Range("J23").Value = 660.2857
$ws.Range("L23").Value = 1980.8571
$ws.Range("N23").Value = -2450.8571
# Row 98: Sweet Kiss of Death
$ws.Range("H98").Value = 450.5
$ws.Range("J98").Value = 483
$ws.Range("L98").Value = 1449
$ws.Range("N98").Value = -4445
# Row 133: Friends Are Food
$ws.Range("H133").Value = 5432.5
$ws.Range("J133").Value = 4999.375
$ws.Range("L133").Value = 14998.125
$ws.Range("N133").Value = -25118.125

# ==== GSM ====
$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers
$ws.Range("H2").Value = 1380.3889
$ws.Range("I2").Value = 97.44444
$ws.Range("J2").Value = 2663.3333
$ws.Range("K2").Value = 97.44444
$ws.Range("L2").Value = 2663.3333
$ws.Range("M2").Value = 15.55556
$ws.Range("N2").Value = -2889.3333
# Row 39: One Man's Trash
$ws.Range("H39").Value = 25506.666
$ws.Range("J39").Value = 25506.666
$ws.Range("L39").Value = 25506.666
$ws.Range("N39").Value = -26570.666
# Row 75: Citizen's Arrest
$ws.Range("H75").Value = 58665.668
$ws.Range("J75").Value = 58665.668
$ws.Range("L75").Value = 58665.668
$ws.Range("N75").Value = -60413.668
# Row 78: Watchers within the Walls (L)
$ws.Range("H78").Value = 58665.668
$ws.Range("J78").Value = 58665.668
$ws.Range("L78").Value = 175997.004
$ws.Range("N78").Value = -184733.004
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 3643.75
$ws.Range("J80").Value = 3982.8
$ws.Range("L80").Value = 3982.8
$ws.Range("N80").Value = -5978.8
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 3643.75
$ws.Range("J83").Value = 3982.8
$ws.Range("L83").Value = 19914
$ws.Range("N83").Value = -29898
# Row 99: Needle in a Hingan Stack
$ws.Range("H99").Value = 14483.167
$ws.Range("I99").Value = 1724.75
$ws.Range("K99").Value = 1724.75
$ws.Range("M99").Value = 521.25
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 4492.8335
$ws.Range("J122").Value = 9514.666999999999
$ws.Range("L122").Value = 28544.001
$ws.Range("N122").Value = -33444.001
# Row 126: Gold Rush Order
$ws.Range("H126").Value = 4333
$ws.Range("I126").Value = 3065.6667
$ws.Range("J126").Value = 4966.6665
$ws.Range("K126").Value = 9197.000100000001
$ws.Range("L126").Value = 14899.9995
$ws.Range("M126").Value = -6727.000100000001
$ws.Range("N126").Value = -19839.9995
# Row 128: To Fight at Her Side
$ws.Range("H128").Value = 37499.332
$ws.Range("J128").Value = 49999.332
$ws.Range("L128").Value = 49999.332
$ws.Range("N128").Value = -59959.332

# ==== LTW ====
$ws = $wb.Worksheets.Item("LTW")
# Row 108: Girding for Glory
$ws.Range("H108").Value = 80620.8
$ws.Range("J108").Value = 80620.8
$ws.Range("L108").Value = 80620.8
$ws.Range("N108").Value = -88300.8
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 3555.0527
$ws.Range("I132").Value = 3344.1177
$ws.Range("K132").Value = 10032.3531
$ws.Range("M132").Value = -7502.3531
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4534.0234
$ws.Range("I136").Value = 4140
$ws.Range("K136").Value = 12420
$ws.Range("M136").Value = -9870

# ==== WVR ====
$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 185299.8
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248
# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 185299.8
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 2199.5417
$ws.Range("I122").Value = 1799.3529
$ws.Range("J122").Value = 3171.4285
$ws.Range("K122").Value = 5398.0587
$ws.Range("L122").Value = 9514.2855
$ws.Range("M122").Value = -2948.0587
$ws.Range("N122").Value = -14414.2855
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 4085.389
$ws.Range("I126").Value = 3971.125
$ws.Range("K126").Value = 11913.375
$ws.Range("M126").Value = -9443.375
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 9096242
$ws.Range("I136").Value = 11148962
$ws.Range("J136").Value = 5627.857
$ws.Range("K136").Value = 33446886
$ws.Range("L136").Value = 16883.571
$ws.Range("M136").Value = -33444336
$ws.Range("N136").Value = -21983.571
